$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 43, shifting existing rows 43:51 down to 44:52
$ws.Rows.Item(43).Insert()

# Copy the formatting/style context of the row that is now 44 (previously row 43)
# for the new row 43's date cell, matching the s="2" style used by column D.
$ws.Cells.Item(43, 4).NumberFormat = $ws.Cells.Item(44, 4).NumberFormat

# Populate the new row 43 with data (same as old row 43 except Fecha and Volumen)
$ws.Cells.Item(43, 1).Value = 4
$ws.Cells.Item(43, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(43, 3).Value = "Los Lagos"
$ws.Cells.Item(43, 4).Value = 44504
$ws.Cells.Item(43, 5).Value = 10
$ws.Cells.Item(43, 6).Value = 100112026
$ws.Cells.Item(43, 7).Value = "Haba"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 70
$ws.Cells.Item(43, 11).Value = 11000
$ws.Cells.Item(43, 12).Value = 11000
$ws.Cells.Item(43, 13).Value = 11000
$ws.Cells.Item(43, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(43, 15).Value = "Región Metropolitana"
$ws.Cells.Item(43, 16).Value = 440
$ws.Cells.Item(43, 17).Value = 25
$ws.Cells.Item(43, 18).Value = "Hortaliza"
